$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  17"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Type-transition cells (text -> number needing explicit NumberFormat) ---
$ws.Cells.Item(14, 4).NumberFormat = "#,##0"
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(14, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(14, 5).Value = -100
$ws.Cells.Item(22, 4).NumberFormat = "#,##0"
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(22, 5).Value = 50
$ws.Cells.Item(30, 3).NumberFormat = "#,##0"
$ws.Cells.Item(30, 3).Value = 1

# --- Type-transition cells (number -> text needing explicit NumberFormat) ---
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "0"

# --- Regular numeric value updates ---
$ws.Cells.Item(14, 10).Value = 13
$ws.Cells.Item(14, 11).Value = -69.230769230769
$ws.Cells.Item(14, 12).Value = -76.470588235294
$ws.Cells.Item(14, 14).Value = -92.156862745098
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 300
$ws.Cells.Item(15, 6).Value = 9
$ws.Cells.Item(15, 7).Value = 11
$ws.Cells.Item(15, 8).Value = -18.181818181818
$ws.Cells.Item(15, 9).Value = 50
$ws.Cells.Item(15, 10).Value = 58
$ws.Cells.Item(15, 11).Value = -13.793103448275
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 21.951219512195
$ws.Cells.Item(15, 14).Value = -54.128440366972
$ws.Cells.Item(16, 3).Value = 24
$ws.Cells.Item(16, 4).Value = 14
$ws.Cells.Item(16, 5).Value = 71.428571428571
$ws.Cells.Item(16, 6).Value = 100
$ws.Cells.Item(16, 7).Value = 94
$ws.Cells.Item(16, 8).Value = 6.382978723404
$ws.Cells.Item(16, 9).Value = 419
$ws.Cells.Item(16, 10).Value = 422
$ws.Cells.Item(16, 11).Value = -0.710900473933
$ws.Cells.Item(16, 12).Value = 43.493150684931
$ws.Cells.Item(16, 13).Value = -38.017751479289
$ws.Cells.Item(16, 14).Value = -82.792607802874
$ws.Cells.Item(17, 3).Value = 51
$ws.Cells.Item(17, 4).Value = 54
$ws.Cells.Item(17, 5).Value = -5.555555555555
$ws.Cells.Item(17, 6).Value = 213
$ws.Cells.Item(17, 7).Value = 185
$ws.Cells.Item(17, 8).Value = 15.135135135135
$ws.Cells.Item(17, 9).Value = 842
$ws.Cells.Item(17, 10).Value = 778
$ws.Cells.Item(17, 11).Value = 8.226221079691
$ws.Cells.Item(17, 12).Value = 24.37223042836
$ws.Cells.Item(17, 13).Value = 67.395626242544
$ws.Cells.Item(17, 14).Value = -21.74721189591
$ws.Cells.Item(18, 3).Value = 23
$ws.Cells.Item(18, 4).Value = 18
$ws.Cells.Item(18, 5).Value = 27.777777777777
$ws.Cells.Item(18, 6).Value = 86
$ws.Cells.Item(18, 7).Value = 70
$ws.Cells.Item(18, 8).Value = 22.857142857142
$ws.Cells.Item(18, 9).Value = 357
$ws.Cells.Item(18, 10).Value = 336
$ws.Cells.Item(18, 11).Value = 6.25
$ws.Cells.Item(18, 12).Value = 27.5
$ws.Cells.Item(18, 13).Value = -41.282894736842
$ws.Cells.Item(19, 3).Value = 64
$ws.Cells.Item(19, 4).Value = 68
$ws.Cells.Item(19, 5).Value = -5.882352941176
$ws.Cells.Item(19, 7).Value = 287
$ws.Cells.Item(19, 8).Value = -12.891986062717
$ws.Cells.Item(19, 9).Value = 1086
$ws.Cells.Item(19, 10).Value = 1206
$ws.Cells.Item(19, 11).Value = -9.950248756218
$ws.Cells.Item(19, 12).Value = 48.158253751705
$ws.Cells.Item(19, 13).Value = 21.070234113712
$ws.Cells.Item(19, 14).Value = -52.514210756449
$ws.Cells.Item(20, 3).Value = 36
$ws.Cells.Item(20, 4).Value = 26
$ws.Cells.Item(20, 5).Value = 38.461538461538
$ws.Cells.Item(20, 6).Value = 124
$ws.Cells.Item(20, 7).Value = 107
$ws.Cells.Item(20, 8).Value = 15.887850467289
$ws.Cells.Item(20, 9).Value = 534
$ws.Cells.Item(20, 10).Value = 555
$ws.Cells.Item(20, 11).Value = -3.783783783783
$ws.Cells.Item(20, 12).Value = 93.478260869565
$ws.Cells.Item(20, 13).Value = 0.754716981132
$ws.Cells.Item(20, 14).Value = -91.570639305445
$ws.Cells.Item(21, 3).Value = 202
$ws.Cells.Item(21, 4).Value = 183
$ws.Cells.Item(21, 5).Value = 10.382513661202
$ws.Cells.Item(21, 6).Value = 784
$ws.Cells.Item(21, 7).Value = 757
$ws.Cells.Item(21, 8).Value = 3.566710700132
$ws.Cells.Item(21, 9).Value = 3292
$ws.Cells.Item(21, 10).Value = 3368
$ws.Cells.Item(21, 11).Value = -2.256532066508
$ws.Cells.Item(21, 12).Value = 41.591397849462
$ws.Cells.Item(21, 13).Value = 0.549786194257
$ws.Cells.Item(21, 14).Value = -78.342105263157
$ws.Cells.Item(22, 3).Value = 3
$ws.Cells.Item(22, 6).Value = 9
$ws.Cells.Item(22, 7).Value = 6
$ws.Cells.Item(22, 8).Value = 50
$ws.Cells.Item(22, 9).Value = 39
$ws.Cells.Item(22, 10).Value = 39
$ws.Cells.Item(22, 12).Value = 77.272727272727
$ws.Cells.Item(22, 13).Value = 2.631578947368
$ws.Cells.Item(23, 3).Value = 4
$ws.Cells.Item(23, 4).Value = 4
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 9).Value = 79
$ws.Cells.Item(23, 10).Value = 72
$ws.Cells.Item(23, 11).Value = 9.722222222222
$ws.Cells.Item(23, 12).Value = 21.538461538461
$ws.Cells.Item(23, 13).Value = 49.056603773584
$ws.Cells.Item(24, 3).Value = 194
$ws.Cells.Item(24, 4).Value = 182
$ws.Cells.Item(24, 5).Value = 6.593406593406
$ws.Cells.Item(24, 6).Value = 637
$ws.Cells.Item(24, 7).Value = 747
$ws.Cells.Item(24, 8).Value = -14.725568942436
$ws.Cells.Item(24, 9).Value = 3032
$ws.Cells.Item(24, 10).Value = 2967
$ws.Cells.Item(24, 11).Value = 2.190765082574
$ws.Cells.Item(24, 12).Value = 41.682242990654
$ws.Cells.Item(24, 13).Value = 59.83131259884
$ws.Cells.Item(25, 3).Value = 76
$ws.Cells.Item(25, 4).Value = 69
$ws.Cells.Item(25, 5).Value = 10.144927536231
$ws.Cells.Item(25, 6).Value = 332
$ws.Cells.Item(25, 7).Value = 277
$ws.Cells.Item(25, 8).Value = 19.85559566787
$ws.Cells.Item(25, 9).Value = 1366
$ws.Cells.Item(25, 10).Value = 1143
$ws.Cells.Item(25, 11).Value = 19.510061242344
$ws.Cells.Item(25, 12).Value = 41.848390446521
$ws.Cells.Item(25, 13).Value = -4.675505931612
$ws.Cells.Item(26, 3).Value = 6
$ws.Cells.Item(26, 4).Value = 3
$ws.Cells.Item(26, 5).Value = 100
$ws.Cells.Item(26, 6).Value = 16
$ws.Cells.Item(26, 7).Value = 24
$ws.Cells.Item(26, 8).Value = -33.333333333333
$ws.Cells.Item(26, 9).Value = 87
$ws.Cells.Item(26, 10).Value = 99
$ws.Cells.Item(26, 11).Value = -12.121212121212
$ws.Cells.Item(26, 12).Value = 6.097560975609
$ws.Cells.Item(27, 3).Value = 13
$ws.Cells.Item(27, 5).Value = 85.714285714285
$ws.Cells.Item(27, 6).Value = 35
$ws.Cells.Item(27, 7).Value = 25
$ws.Cells.Item(27, 8).Value = 40
$ws.Cells.Item(27, 9).Value = 136
$ws.Cells.Item(27, 10).Value = 129
$ws.Cells.Item(27, 11).Value = 5.426356589147
$ws.Cells.Item(27, 12).Value = 19.298245614035
$ws.Cells.Item(28, 4).Value = 8
$ws.Cells.Item(28, 5).Value = -87.5
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(28, 7).Value = 15
$ws.Cells.Item(28, 8).Value = -66.666666666666
$ws.Cells.Item(28, 9).Value = 35
$ws.Cells.Item(28, 10).Value = 50
$ws.Cells.Item(28, 11).Value = -30
$ws.Cells.Item(28, 12).Value = -12.5
$ws.Cells.Item(28, 13).Value = -40.677966101694
$ws.Cells.Item(28, 14).Value = -79.411764705882
$ws.Cells.Item(29, 4).Value = 6
$ws.Cells.Item(29, 5).Value = -83.333333333333
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(29, 7).Value = 12
$ws.Cells.Item(29, 8).Value = -58.333333333333
$ws.Cells.Item(29, 9).Value = 26
$ws.Cells.Item(29, 10).Value = 41
$ws.Cells.Item(29, 11).Value = -36.585365853658
$ws.Cells.Item(29, 12).Value = -27.777777777777
$ws.Cells.Item(29, 13).Value = -46.938775510204
$ws.Cells.Item(29, 14).Value = -83.006535947712
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 8
$ws.Cells.Item(30, 7).Value = 8
$ws.Cells.Item(30, 9).Value = 17
$ws.Cells.Item(30, 10).Value = 13
$ws.Cells.Item(30, 11).Value = 30.76923076923
$ws.Cells.Item(30, 12).Value = 750

Write-Host "done"